$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet, positioned between "2021-Q4" and
#    "总计". We get there by duplicating the existing "总计" sheet
#    (so the duplicate inherits its header/index-column styling) and
#    then overwriting its content with the 2022-Q1 fund holdings.
# ------------------------------------------------------------------
$totals0 = $wb.Worksheets.Item("总计")
$totals0.Copy($totals0)
$q1 = $wb.Worksheets.Item("总计 (2)")
$q1.Name = "2022-Q1"
# Re-resolve the original "总计" sheet by name - the handle captured
# before the Copy() call above tracks the sheet *position*, which the
# new copy now occupies, so it must be looked up again post-copy.
$totals = $wb.Worksheets.Item("总计")

# Headers (row 1) - same column layout/style as the "总计" header row,
# extended out to column H to match the "2021-Q4" sheet's layout.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data cells in columns B-G are stored as text (verbatim strings from
# the source report, including leading zeros in fund codes), so force
# a text number format before writing the values.
$q1.Range("B2:G3").NumberFormat = "@"

# Row 2 - 009837 / 华夏磐锐一年定期开放混合A
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "009837"
$q1.Range("C2").Value = "华夏磐锐一年定期开放混合A"
$q1.Range("D2").Value = "16.45"
$q1.Range("E2").Value = "79.44"
$q1.Range("F2").Value = "4.04"
$q1.Range("G2").Value = "0.6646"
$q1.Range("H2").Value = 3

# Row 3 - 009838 / 华夏磐锐一年定期开放混合C
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "009838"
$q1.Range("C3").Value = "华夏磐锐一年定期开放混合C"
$q1.Range("D3").Value = "0.44"
$q1.Range("E3").Value = "79.44"
$q1.Range("F3").Value = "4.04"
$q1.Range("G3").Value = "0.0178"
$q1.Range("H3").Value = 3

# The "@" number format above leaves a stray cell style behind - strip
# it back off (columns B-G in the data rows carry no explicit style in
# the source workbook) by pasting in a pristine, never-formatted cell's
# format.
$q1.Range("Z1").Copy()
$q1.Range("B2:G3").PasteSpecial(-4122)

# The copied "总计" sheet's A2 cell already carries the bold "index
# column" style - replicate it onto A3 too.
$q1.Range("A2").Copy()
$q1.Range("A3").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Update the "总计" sheet: add a 2022-Q1 total row above the
#    existing 2021-Q4 row (which shifts down to row 3).
# ------------------------------------------------------------------
# Move the existing 2021-Q4 totals down to row 3 first.
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2021-Q4"
$totals.Range("C3").Value = 1
$totals.Range("D3").Value = 0

# Carry the bold index-column style from A2 onto the relocated A3.
$totals.Range("A2").Copy()
$totals.Range("A3").PasteSpecial(-4122)

# Now overwrite row 2 with the new 2022-Q1 totals.
$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 0.68
